$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: all columns A-L get 149.07119849998597, except D2 which is 4x that value
$ws.Range("A2:L2").Value = 149.07119849998597
$ws.Range("D2").Value = 596.28479399994387

# Row 3: all columns A-L get 133.33333333333331, except D3 which is 4x that value
$ws.Range("A3:L3").Value = 133.33333333333331
$ws.Range("D3").Value = 533.33333333333326
